# edit.ps1 -- applies the "BBI-23-6" revision described in the diff:
#   - fills in several previously-blank score cells on "Лист1"
#     (H2, H3, I12, E17, G17, H18, H19, H20, H23)
#   - moves the active-cell selection from G19 to H19
#   - wipes the scratch/demo area below the main table (rows 30-35):
#       * row 30 (the little 0..7 demo row) is cleared entirely
#       * row 31 keeps only its four bordered placeholder cells (B/D/E/F), now blank
#       * the "gap/i/j/temp" mini table in A32:B35 has its labels and
#         values removed (left blank, same as the rest of those rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

# --- fill in previously empty grade inputs -----------------------------
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I12").Value = 3
$ws.Range("E17").Value = -2
$ws.Range("G17").Value = -1
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H23").Value = 0

# --- clear the scratch area below the table -----------------------------
# Row 30 (0,1,2,3,4,5,6,7 demo values) -> fully cleared
$ws.Rows(30).ClearContents()

# Row 31 keeps B31/D31/E31/F31 (blank, bordered) but loses C31/G31/H31/I31
$ws.Range("C31").ClearContents()
$ws.Range("G31:I31").ClearContents()

# gap/i/j/temp mini table -> labels and values removed
$ws.Range("A32:B35").ClearContents()

# --- move the selection to match the saved view -------------------------
[void]$ws.Range("H19").Select()
